$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value added
$ws1.Range("B9").Value = "Alvearie Team"

# Replace the first "Contact" row with "Jurisdiction"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Remove the duplicate "Contact" row entirely (row shifts everything below up by one)
$ws1.Rows.Item(11).Delete()

# --- Elements sheet ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root Extension element: Short/Definition now reflect the profile title/description
$ws2.Range("K2").Value = "Disposition Reason"
$ws2.Range("L2").Value = "Customer-specific code for the disposition reason, as related to how the claim was paid"
